# Update "想去人数" (F column) values on several sheets, as generated at
# commit 456a3b4 for gh-pages output.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 149
$ws1.Range("F3").Value = 1344
$ws1.Range("F4").Value = 1152
$ws1.Range("F5").Value = 1044
$ws1.Range("F6").Value = 1829
$ws1.Range("F7").Value = 582
$ws1.Range("F8").Value = 1217
$ws1.Range("F12").Value = 308
$ws1.Range("F13").Value = 85
$ws1.Range("F14").Value = 92
$ws1.Range("F15").Value = 718
$ws1.Range("F16").Value = 188
$ws1.Range("F17").Value = 110
$ws1.Range("F18").Value = 30
$ws1.Range("F21").Value = 168
$ws1.Range("F22").Value = 681
$ws1.Range("F23").Value = 48
$ws1.Range("F25").Value = 168
$ws1.Range("F27").Value = 884
$ws1.Range("F28").Value = 325
$ws1.Range("F29").Value = 167
$ws1.Range("F30").Value = 49
$ws1.Range("F31").Value = 283
$ws1.Range("F34").Value = 410

# Sheet "演出" (sheet2)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F6").Value = 26
$ws2.Range("F10").Value = 621

# Sheet "本地生活" (sheet3)
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 315

# Sheet "全部类型" (sheet4)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 315
$ws4.Range("F3").Value = 149
$ws4.Range("F4").Value = 1343
$ws4.Range("F5").Value = 1152
$ws4.Range("F6").Value = 1044
$ws4.Range("F7").Value = 1829
$ws4.Range("F8").Value = 582
$ws4.Range("F9").Value = 1217
$ws4.Range("F14").Value = 308
$ws4.Range("F15").Value = 85
$ws4.Range("F16").Value = 92
$ws4.Range("F17").Value = 718
$ws4.Range("F18").Value = 188
$ws4.Range("F19").Value = 110
$ws4.Range("F21").Value = 30
$ws4.Range("F26").Value = 26
$ws4.Range("F29").Value = 168
$ws4.Range("F30").Value = 681
$ws4.Range("F31").Value = 48
$ws4.Range("F33").Value = 168
$ws4.Range("F35").Value = 884
$ws4.Range("F36").Value = 325
$ws4.Range("F39").Value = 167
$ws4.Range("F40").Value = 49
$ws4.Range("F41").Value = 283
$ws4.Range("F42").Value = 621
$ws4.Range("F48").Value = 410
